# Updated cryptos list (price + 1h volume change) for the "cryptos" sheet.
# Values in column D are written with a leading apostrophe so Excel keeps
# them as text (preserving formats like "28.099.37", "0.00001108", "5.010")
# instead of silently re-parsing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.099.37"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "'1.873.20"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'313.16"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D7").Value = "'0.5135"
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("D8").Value = "'0.3891"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").Value = "'0.08374"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").Value = "'41.69"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "'6.197"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").Value = "'20.59"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "'1.872.54"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").Value = "'7.288"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "'0.00001108"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").Value = "'90.89"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "'0.06654"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "'17.70"
$ws.Range("E20").Value = "  -1.91%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "'6.028"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "'28.136.79"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'11.13"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("D25").Value = "'2.251"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").Value = "'2.080.47"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").Value = "'2.472"
$ws.Range("E27").Value = "  -4.28%  "
$ws.Range("D28").Value = "'158.14"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").Value = "'20.66"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "'125.62"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("D33").Value = "'5.898"
$ws.Range("E33").Value = "  +4.71%  "
$ws.Range("D34").Value = "'3.599"
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("D35").Value = "'9.689"
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").Value = "'0.06535"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").Value = "'0.2185"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("D40").Value = "'0.6501"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("D41").Value = "'5.010"
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("D43").Value = "'11.31"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").Value = "'0.6095"
$ws.Range("E44").Value = "  -1.53%  "
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("D46").Value = "'1.280"
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("D47").Value = "'3.673"
$ws.Range("D48").Value = "'2.009"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("D49").Value = "'1.216"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "'121.31"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'77.78"
$ws.Range("E51").Value = "  -3.32%  "
